$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B2 mark from 7.5 to 6.5
$ws.Range("B2").Value = 6.5

# Replace row 3 (was Winter/2019, 10, IP-31) with the new Winter/2020, 5, KS-41 row
$ws.Range("A3").Value = "Winter/2020"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "KS-41"

# Delete old rows 4 and 5 entirely (they are no longer part of the table)
$ws.Range("A4:C5").Delete()
